$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 805
    "C2" = 871
    "D2" = 812
    "E2" = 811
    "F2" = 320
    "G2" = 391
    "H2" = 280
    "I2" = 1269
    "J2" = 1269
    "K2" = 1269
    "L2" = 1269
    "M2" = 1269
    "N2" = 1269
    "O2" = 1269
    "P2" = 1269
    "Q2" = 1269
    "R2" = 1269
    "S2" = 1269
    "T2" = 1269
    "U2" = 1269
    "V2" = 1269
    "W2" = 1269
    "X2" = 1269
    "Y2" = 1269
    "Z2" = 1269
    "AA2" = 1269
    "AB2" = 1269
    "AC2" = 1269
    "AD2" = 1269
    "AE2" = 1269
    "AF2" = 1269
    "AG2" = 1269
    "AH2" = 1269
    "AI2" = 1269
    "AJ2" = 1269
    "AK2" = 1269
    "AL2" = 1269
    "AM2" = 1269
    "AN2" = 1269
    "AO2" = 1269
    "AP2" = 1269
    "AQ2" = 1269
    "AR2" = 1269
    "AS2" = 1269
    "AT2" = 1269
    "AU2" = 1269
    "AV2" = 1269
    "AW2" = 1269
    "AX2" = 1269
    "AY2" = 1269
    "AZ2" = 1269
    "BA2" = 1269
    "BB2" = 1269
    "BC2" = 1269
    "BD2" = 1269
    "BE2" = 1269
    "BF2" = 1269
    "BG2" = 1269
    "B3" = 5.472894409937888
    "C3" = 625.4120283199388
    "D3" = 676.4246408045977
    "E3" = 1075.919340320592
    "F3" = 1542.71484375
    "G3" = 94.93711338448422
    "H3" = 10.80166369047619
    "I3" = 7.705899159364328
    "J3" = 1.757968630260047
    "K3" = 2.531229395744681
    "L3" = 0.8972608353033883
    "M3" = 1.969253821907014
    "N3" = 0.7897307342789599
    "O3" = 10.15173316469661
    "P3" = 0.7276483769897557
    "Q3" = 1.721290255476753
    "R3" = 0.8200179100866825
    "S3" = 3.374662150985028
    "T3" = 0.4095350669818755
    "U3" = 0.03393223010244287
    "V3" = 2.067297084318361
    "W3" = 14.30319076246388
    "X3" = 15.55751966485421
    "Y3" = 4.066851112923562
    "Z3" = 6.022582594562648
    "AA3" = 5.882695035460993
    "AB3" = 1.589913317572892
    "AC3" = 0.7491725768321513
    "AD3" = 3.545784081954295
    "AE3" = 0.0425531914893617
    "AF3" = 1.124142927265563
    "AG3" = 0.1245547675334909
    "AH3" = 0.007092198581560284
    "AI3" = 3.776674546887313
    "AJ3" = 0.6028368794326241
    "AK3" = 0.3685579196217494
    "AL3" = 0.1601457840819543
    "AM3" = 0.1599684791174153
    "AN3" = 0.975177304964539
    "AO3" = 0.1646966115051221
    "AP3" = 2.061229314420804
    "AQ3" = 0.1379038613081166
    "AR3" = 0.8490640630417653
    "AS3" = 0.1533718088258472
    "AT3" = 0.1952583845547675
    "AU3" = 0.2736012608353034
    "AV3" = 1.007662577462569
    "AW3" = 0.3386190398739165
    "AX3" = 0.1007722616233254
    "AY3" = 0.07313632781717887
    "AZ3" = 0.07486209613869188
    "BA3" = 0.1095350669818755
    "BB3" = 0.1107171000788022
    "BC3" = 0.001576044129235619
    "BD3" = 0.2442868400315209
    "BE3" = 0.08668242710795902
    "B4" = 5.430753289239733
    "C4" = 172.2104111356757
    "D4" = 176.7289916310277
    "E4" = 265.9006960655318
    "F4" = 502.9354187947058
    "G4" = 52.99942205420455
    "H4" = 12.93810686313891
    "I4" = 14.43346693300467
    "J4" = 5.287464501756308
    "K4" = 6.754250447766093
    "L4" = 3.697799905897668
    "M4" = 5.69952924215833
    "N4" = 2.840631030449631
    "O4" = 21.58336133597242
    "P4" = 2.984730810999506
    "Q4" = 5.243466376814968
    "R4" = 3.792461522792751
    "S4" = 8.241942753428274
    "T4" = 5.093805744266537
    "U4" = 0.6672454311834292
    "V4" = 11.18125800714043
    "W4" = 16.48406269562457
    "X4" = 22.78300302513044
    "Y4" = 9.169522313736033
    "Z4" = 7.929484486384751
    "AA4" = 16.57990182793024
    "AB4" = 6.414745206555939
    "AC4" = 5.041524075995074
    "AD4" = 13.90404858274003
    "AE4" = 0.4279307261182981
    "AF4" = 5.25364391211117
    "AG4" = 0.834141577602706
    "AH4" = 0.204322706292276
    "AI4" = 13.38679694571391
    "AJ4" = 3.430624947347729
    "AK4" = 4.255408088259993
    "AL4" = 2.359522955787477
    "AM4" = 2.368718057430808
    "AN4" = 7.678756598970639
    "AO4" = 2.56409606879239
    "AP4" = 11.12290715019279
    "AQ4" = 2.14366994645412
    "AR4" = 4.041652497877855
    "AS4" = 1.509991139991264
    "AT4" = 2.874846213071987
    "AU4" = 3.971889160505682
    "AV4" = 4.495458317642192
    "AW4" = 3.098775666522288
    "AX4" = 1.891608824712009
    "AY4" = 0.9384063935607984
    "AZ4" = 1.908367780967372
    "BA4" = 2.273800620802517
    "BB4" = 2.06646161398685
    "BC4" = 0.05614346140443382
    "BD4" = 3.510267845432892
    "BE4" = 1.781535361340087
    "F5" = 140.5
    "C6" = 426
    "D6" = 477
    "E6" = 792.2
    "F6" = 1304.75
    "G6" = 61.18
    "H6" = 1.7
    "C7" = 668.4
    "E7" = 1145
    "F7" = 1581.5
    "G7" = 88.7
    "H7" = 5.6
    "W7" = 9.4
    "C8" = 740
    "D8" = 796.8499999999999
    "E8" = 1267.5
    "F8" = 1797.625
    "G8" = 108.4
    "H8" = 14.625
    "I8" = 10
    "O8" = 3
    "W8" = 24.75
    "Y8" = 0
    "AS9" = 19.2
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
